$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the received-ID value for Henri LESAGE (row 2), keeping the cell itself intact.
$ws.Range("D2").ClearContents()

# Update "Total sans reçus" count to 0.
$ws.Range("B6").Value = 0
